$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  @(2000,140,5,20),
  @(2000,140,5,20),
  @(2000,140,5,20),
  @(1000,140,5,20),
  @(1000,140,5,20),
  @(1000,140,5,20),
  @(1000,140,5,20),
  @(1000,140,5,20),
  @(1000,140,5,20)
)

$startRow = 17
for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $startRow + $i
  $rowVals = $values[$i]
  $ws.Cells.Item($row, 1).Value = $rowVals[0]
  $ws.Cells.Item($row, 2).Value = $rowVals[1]
  $ws.Cells.Item($row, 3).Value = $rowVals[2]
  $ws.Cells.Item($row, 4).Value = $rowVals[3]
}
